$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "64.334.69"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "3.094.80"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue $ws.Range("D5") "559.97"
$ws.Range("E5").Value = "  +1.88%  "

Set-TextValue $ws.Range("D6") "144.24"
$ws.Range("E6").Value = "  +2.89%  "

Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.091.20"
$ws.Range("E8").Value = "  +1.03%  "

Set-TextValue $ws.Range("D9") "0.507"
$ws.Range("E9").Value = "  +1.05%  "

Set-TextValue $ws.Range("D10") "0.154"
$ws.Range("E10").Value = "  +1.34%  "

Set-TextValue $ws.Range("D11") "6.17"
$ws.Range("E11").Value = "  -5.86%  "

Set-TextValue $ws.Range("D12") "0.472"
$ws.Range("E12").Value = "  +3.75%  "

Set-TextValue $ws.Range("D13") "0.0000229"
$ws.Range("E13").Value = "  +0.45%  "

Set-TextValue $ws.Range("D14") "35.17"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").Value = "3.589.10"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").Value = "64.319.67"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "3.091.09"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("E18").Value = "  +1.23%  "

Set-TextValue $ws.Range("D19") "6.76"
$ws.Range("E19").Value = "  -0.23%  "

Set-TextValue $ws.Range("D20") "486.89"
$ws.Range("E20").Value = "  +0.41%  "

Set-TextValue $ws.Range("D21") "13.99"
$ws.Range("E21").Value = "  +1.58%  "

Set-TextValue $ws.Range("D22") "0.676"
$ws.Range("E22").Value = "  -0.08%  "

Set-TextValue $ws.Range("D23") "7.58"
$ws.Range("E23").Value = "  +3.89%  "

Set-TextValue $ws.Range("D24") "14.19"
$ws.Range("E24").Value = "  +12.32%  "

Set-TextValue $ws.Range("D25") "81.33"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  +0.13%  "

Set-TextValue $ws.Range("D27") "2.80"
$ws.Range("E27").Value = "  +1.34%  "

Set-TextValue $ws.Range("D28") "8.03"
$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("E29").Value = "  +2.77%  "

Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.33%  "

Set-TextValue $ws.Range("D31") "26.44"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("E33").Value = "  +1.29%  "

Set-TextValue $ws.Range("D34") "5.60"
$ws.Range("E34").Value = "  -1.93%  "

Set-TextValue $ws.Range("D35") "6.25"
$ws.Range("E35").Value = "  +4.29%  "

Set-TextValue $ws.Range("D36") "55.87"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D37") "3.02"
$ws.Range("E37").Value = "  +18.03%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D38") "452.85"
$ws.Range("E38").Value = "  -3.05%  "

Set-TextValue $ws.Range("D39") "0.0409"
$ws.Range("E39").Value = "  +2.79%  "

$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("D41").Value = "2.972.48"
$ws.Range("E41").Value = "  -2.86%  "

Set-TextValue $ws.Range("D42") "8.24"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("E43").Value = "  -5.81%  "

Set-TextValue $ws.Range("D44") "28.12"
$ws.Range("E44").Value = "  -0.44%  "

Set-TextValue $ws.Range("D45") "0.262"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("E46").Value = "  -0.02%  "

Set-TextValue $ws.Range("D47") "2.14"
$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("E48").Value = "  +1.97%  "

Set-TextValue $ws.Range("D49") "118.98"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("E50").Value = "  +0.67%  "

Set-TextValue $ws.Range("D51") "2.09"
$ws.Range("E51").Value = "  +0.18%  "
